$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new row above row 1; this shifts all existing rows (1-47) down to (2-48)
#    and keeps their original formatting/values intact.
$ws.Rows("1:1").Insert()

# 2. Merge the new header row's A1:E1 range for the legend text.
$ws.Range("A1:E1").Merge()

# 3. Build the rich-text legend value: bold title followed by a normal-weight
#    explanatory paragraph.
$title = "Significantly co-expressed GO terms"
$body = "`nCo-expression was measured among genes within each GO term that had co-expression data in each network using both density (Eq. 1) and locality (Eq. 2). Significance of co-expression metrics was assessed by comparing values to 1,000 random gene sets of the same size.`n"
$ws.Range("A1").Value = $title + $body

$titleRun = $ws.Range("A1").Characters(1, $title.Length)
$titleRun.Font.Name = "Calibri"
$titleRun.Font.Size = 11
$titleRun.Font.Bold = $true

$bodyRun = $ws.Range("A1").Characters($title.Length + 1, $body.Length)
$bodyRun.Font.Name = "Calibri"
$bodyRun.Font.Size = 11
$bodyRun.Font.Bold = $false

# 4. Apply the "legend" border/style to the rest of the header row by copying the
#    existing bottom+left border style used elsewhere in the sheet (now located at
#    A6 after the row insert) and stripping its left edge, producing a bottom-only
#    border that matches the rest of the sheet's look.
$ws.Range("A6").Copy()
$ws.Range("B1:E1").PasteSpecial(-4122)
$ws.Range("B1:E1").Borders(7).LineStyle = -4142

$ws.Range("A6").Copy()
$ws.Range("A1").PasteSpecial(-4122)
$ws.Range("A1").Borders(7).LineStyle = -4142
$ws.Range("A1").WrapText = $true

$excel.CutCopyMode = $false

# 5. Give the legend row extra height so the wrapped paragraph is fully visible.
$ws.Rows("1:1").RowHeight = 80

# 6. Freeze the legend row so the table header/body stay visible while scrolling,
#    then restore the cursor to the cell that was selected after the edit.
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("E18").Select()
